$d = $word.ActiveDocument

# Locate the paragraph containing the "git:remote -a siraj-ass5-midterm" command
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*git:remote*") {
        $target = $p
    }
}

$r = $target.Range
$find = $r.Find
$find.Text = "siraj-ass5-midterm"
$find.Replacement.Text = "siraj-ass10-final"
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, "siraj-ass10-final", 2)
